$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 5

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 3
